# Applies the "Add files via upload" commit: the survey export was
# re-run later the same day (new two responses came in), producing an
# updated raw-data sheet:
#   - sheet tab renamed to reflect the newer export timestamp
#   - two additional respondent rows appended (row 12 and row 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (new export timestamp: 17:15 -> 18:20)
$ws.Name = "原始数据_202410221820_0"

# New survey response row 12 (answer #11)
$ws.Range("A12").Value = 11
$ws.Range("C12").Value = "填空1:男|填空2:22"
$ws.Range("D12").Value = "B.0-5次"
$ws.Range("E12").Value = "A.0-1小时"
$ws.Range("F12").Value = "C.睡眠改善"
$ws.Range("G12").Value = "A.工作时间过长，没时间锻炼"
$ws.Range("H12").Value = "A.智能手环/手表|B.健康APP（如Keep、MyFitnessPal）"
$ws.Range("I12").Value = "A.是"
$ws.Range("J12").Value = "B.个性化运动计划"
$ws.Range("K12").Value = "B.一般"
$ws.Range("L12").Value = "无"
$ws.Range("M12").Value = "A.无所谓"
$ws.Range("N12").Value = "A.是"
$ws.Range("O12").Value = "C.获取专业指导"
$ws.Range("P12").Value = "C.团体讨论"
$ws.Range("Q12").Value = "B.饮食建议"
$ws.Range("R12").Value = "2024-10-22 17:41:16"
$ws.Range("S12").Value = "2024-10-22 17:44:50"
$ws.Range("T12").Value = "3分34秒"
$ws.Range("U12").Value = "山东省"
$ws.Range("V12").Value = "滨州市"
$ws.Range("W12").Value = "39.144.109.96"
$ws.Range("X12").Value = "Chrome 86.0.4240.99"
$ws.Range("Y12").Value = "Android Linux 10"

# New survey response row 13 (answer #12)
$ws.Range("A13").Value = 12
$ws.Range("C13").Value = "填空1:男|填空2:20"
$ws.Range("D13").Value = "A.0次"
$ws.Range("E13").Value = "A.0-1小时"
$ws.Range("F13").Value = "A.体重管理|B.心理健康|C.睡眠改善"
$ws.Range("G13").Value = "A.工作时间过长，没时间锻炼|B.遗传因素|C.工作性质损伤身体"
$ws.Range("H13").Value = "A.智能手环/手表|B.健康APP（如Keep、MyFitnessPal）|C.远程医疗服务|D.家庭健康监测设备"
$ws.Range("I13").Value = "A.是"
$ws.Range("J13").Value = "A.基于个人数据的定制化饮食建议|B.个性化运动计划"
$ws.Range("K13").Value = "A.没有"
$ws.Range("L13").Value = "可以更流畅一点"
$ws.Range("M13").Value = "C.很重要"
$ws.Range("N13").Value = "A.是"
$ws.Range("O13").Value = "A.分享健康经验|B.寻找健康伙伴"
$ws.Range("P13").Value = "A.评论区留言|B.私信聊天"
$ws.Range("Q13").Value = "A.健身技巧分享|B.饮食建议"
$ws.Range("R13").Value = "2024-10-22 18:10:58"
$ws.Range("S13").Value = "2024-10-22 18:11:36"
$ws.Range("T13").Value = "0分38秒"
$ws.Range("U13").Value = "山东省"
$ws.Range("V13").Value = "青岛市"
$ws.Range("W13").Value = "112.224.155.13"
$ws.Range("X13").Value = "Unknown Browser"
$ws.Range("Y13").Value = "iPhone iOS 18.0.1"
